$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace "Vurdere kriterier for de ulike skytjenestene" with the new task text
# across all rows that currently reference it (B5:B7).
$ws.Range("B5").Value = "Lese om og sammenligne skytjenester"
$ws.Range("B6").Value = "Lese om og sammenligne skytjenester"
$ws.Range("B7").Value = "Lese om og sammenligne skytjenester"

# Shift the remaining task labels down by one row: the "Sette opp sensor
# nettverk m/Raspberry Pi" task now also covers B11, "Implementere
# databaseløsning m/webløsning" now also covers B15, and the
# "Ferdigstilling av rapport " task starts one row later (B16:B19 only).
$ws.Range("B11").Value = "Sette opp sensor nettverk m/Raspberry Pi"
$ws.Range("B15").Value = "Implementere databaseløsning m/webløsning"
$ws.Range("B16").Value = "Ferdigstilling av rapport "
$ws.Range("B17").Value = "Ferdigstilling av rapport "
$ws.Range("B18").Value = "Ferdigstilling av rapport "
$ws.Range("B19").Value = "Ferdigstilling av rapport "

# Update the view state: scroll so row 15 is at the top and move the
# active cell/selection to C12.
$ws.Range("C12").Select()
$excel.ActiveWindow.ScrollRow = 15
